# Update figures for report
# Sheet "summary2" rows 2-3 (EEC_PR.q75 / EEC_high) get refreshed numbers and
# the two row labels trade places: the row that used to hold "EEC_PR.q75"
# now reports the new "EEC_high" figures, and the row that used to hold
# "EEC_high" now reports the (previously "EEC_PR.q75") figures that used to
# live in row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("summary2")

# --- Row labels ---
$ws.Range("B2").Value = "EEC_high"
$ws.Range("B3").Value = "EEC_PR.q75"

# --- Row 2 ("EEC_high") : new figures ---
$ws.Range("C2:G2").Value = 0.0
$ws.Range("H2").Value = 4.747626186906547
$ws.Range("I2").Value = 4.897551224387806
$ws.Range("J2").Value = 5.347326336831585
$ws.Range("K2").Value = 11.994002998500749
$ws.Range("L2").Value = 13.943028485757122
$ws.Range("M2:Q2").Value = 100.0

# --- Row 3 ("EEC_PR.q75") : figures formerly shown in row 2 ---
$ws.Range("C3:G3").Value = 9.507898858448504
$ws.Range("H3").Value = 15.876655277095423
$ws.Range("I3").Value = 16.20603702879181
$ws.Range("J3").Value = 16.696975459214272
$ws.Range("K3").Value = 15.798828077014075
$ws.Range("L3").Value = 19.327900164312137
$ws.Range("M3").Value = 35.58942756269336
$ws.Range("N3").Value = 34.70299684258166
$ws.Range("O3").Value = 42.123551802056056
$ws.Range("P3").Value = 44.12255125271142
$ws.Range("Q3").Value = 53.786050736833246
